# Update "想去人数" (column F) values across the four worksheets to reflect
# newly scraped counts (gh-pages output regenerated at commit 456a3b4).

$wb = $excel.ActiveWorkbook

function Set-FValues {
    param([string]$SheetName, [hashtable]$Updates)
    $ws = $wb.Worksheets.Item($SheetName)
    foreach ($row in $Updates.Keys) {
        $ws.Range("F$row").Value = $Updates[$row]
    }
}

# Sheet: 展览
Set-FValues "展览" @{
    2  = 892
    4  = 4256
    6  = 413
    7  = 3386
    8  = 950
    11 = 268
    12 = 2289
    13 = 1242
    14 = 24
    15 = 1981
    16 = 488
    17 = 237
    18 = 48
    19 = 9422
    20 = 5880
    21 = 374
    22 = 197
    23 = 796
    24 = 95
    25 = 818
    28 = 960
    29 = 442
    30 = 90
    32 = 197
    33 = 4762
    35 = 981
    36 = 118
    37 = 445
}

# Sheet: 演出
Set-FValues "演出" @{
    15 = 3510
    16 = 74
}

# Sheet: 本地生活
Set-FValues "本地生活" @{
    2 = 8601
    3 = 406
    4 = 1473
}

# Sheet: 全部类型
Set-FValues "全部类型" @{
    2  = 8601
    3  = 892
    4  = 406
    5  = 1473
    7  = 4256
    9  = 413
    10 = 3386
    11 = 950
    14 = 268
    15 = 2289
    20 = 1242
    22 = 24
    24 = 488
    25 = 237
    26 = 48
    27 = 9423
    28 = 3510
    29 = 74
    30 = 374
    31 = 197
    32 = 796
    33 = 95
    34 = 818
    36 = 960
    37 = 442
    38 = 90
    41 = 197
    42 = 4762
    44 = 981
    45 = 445
}
